$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.644726333333334
$ws.Range("H2").Value = 4.934179
$ws.Range("I2").Value = 0.03084360558270512
$ws.Range("J2").Value = 0.03084360558270512
$ws.Range("M2").Value = 0.8685706666666667
$ws.Range("N2").Value = 2.605712
$ws.Range("O2").Value = 0.2707495698024546
$ws.Range("P2").Value = 0.2707495698024546
$ws.Range("Q2").Value = 1.428561047827556
$ws.Range("R2").Value = 12.857049430448
$ws.Range("S2").Value = 0.008350892942673998
$ws.Range("T2").Value = 0.008350892942673998

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.644726333333334
$ws.Range("H3").Value = 4.934179
$ws.Range("I3").Value = 0.03084360558270512
$ws.Range("J3").Value = 0.03084360558270512
$ws.Range("O3").Value = 0.4617398237148598
$ws.Range("P3").Value = 0.4617398237148598
$ws.Range("Q3").Value = 2.436286516987222
$ws.Range("R3").Value = 21.926578652885
$ws.Range("S3").Value = 0.01424172100448893
$ws.Range("T3").Value = 0.01424172100448893

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.644726333333334
$ws.Range("H4").Value = 4.934179
$ws.Range("I4").Value = 0.03084360558270512
$ws.Range("J4").Value = 0.03084360558270512
$ws.Range("O4").Value = 0.2675106064826855
$ws.Range("P4").Value = 0.2675106064826855
$ws.Range("Q4").Value = 1.41147124474
$ws.Range("R4").Value = 12.70324120266
$ws.Range("S4").Value = 0.008250991635542193
$ws.Range("T4").Value = 0.008250991635542193

$ws.Range("I5").Value = 0.828024694817689
$ws.Range("J5").Value = 0.828024694817689
$ws.Range("M5").Value = 0.8685706666666667
$ws.Range("N5").Value = 2.605712
$ws.Range("O5").Value = 0.2707495698024546
$ws.Range("P5").Value = 0.2707495698024546
$ws.Range("Q5").Value = 38.35102295300801
$ws.Range("R5").Value = 345.159206577072
$ws.Range("S5").Value = 0.2241873299076981
$ws.Range("T5").Value = 0.2241873299076981

$ws.Range("I6").Value = 0.828024694817689
$ws.Range("J6").Value = 0.828024694817689
$ws.Range("O6").Value = 0.4617398237148598
$ws.Range("P6").Value = 0.4617398237148598
$ws.Range("S6").Value = 0.3823319766166703
$ws.Range("T6").Value = 0.3823319766166703

$ws.Range("I7").Value = 0.828024694817689
$ws.Range("J7").Value = 0.828024694817689
$ws.Range("O7").Value = 0.2675106064826855
$ws.Range("P7").Value = 0.2675106064826855
$ws.Range("S7").Value = 0.2215053882933206
$ws.Range("T7").Value = 0.2215053882933206

$ws.Range("I8").Value = 0.1411316995996059
$ws.Range("J8").Value = 0.1411316995996059
$ws.Range("M8").Value = 0.8685706666666667
$ws.Range("N8").Value = 2.605712
$ws.Range("O8").Value = 0.2707495698024546
$ws.Range("P8").Value = 0.2707495698024546
$ws.Range("Q8").Value = 6.536695203194667
$ws.Range("R8").Value = 58.83025682875201
$ws.Range("S8").Value = 0.03821134695208255
$ws.Range("T8").Value = 0.03821134695208255

$ws.Range("I9").Value = 0.1411316995996059
$ws.Range("J9").Value = 0.1411316995996059
$ws.Range("O9").Value = 0.4617398237148598
$ws.Range("P9").Value = 0.4617398237148598
$ws.Range("S9").Value = 0.06516612609370058
$ws.Range("T9").Value = 0.06516612609370058

$ws.Range("I10").Value = 0.1411316995996059
$ws.Range("J10").Value = 0.1411316995996059
$ws.Range("O10").Value = 0.2675106064826855
$ws.Range("P10").Value = 0.2675106064826855
$ws.Range("S10").Value = 0.03775422655382276
$ws.Range("T10").Value = 0.03775422655382276

